# "Generate Report for Handback"
#
# Updates the localization-status workbook to reflect that both rows
# (the 42a1ff46... and 63a46863... source files) have been handed back
# and are in sync with en-US:
#   - Status column ("Ready for handoff" -> "Handed back: in sync with en-US")
#     on both the zh-cn and de-de sheets.
#   - "Latest Handback DateTime" (column H) populated with a real timestamp
#     instead of the zero-date placeholder.
#   - Two new populated columns: "Latest Target File" (F) and
#     "Latest Handback File" (G), each holding a hyperlink mirroring the
#     existing Source File Name / Latest Handoff File hyperlinks.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

function Update-LangSheet($Workbook, $SheetName, $Lang, $HandbackDateTime, $HandoffCommit, $StatusText) {

    $ws = $Workbook.Worksheets.Item($SheetName)

    $rows = @(
        @{ Row = 2; Guid = "42a1ff46-fda8-4478-9547-2b3f38447945" },
        @{ Row = 3; Guid = "63a46863-c1fb-4ae1-b814-cead6937743c" }
    )

    foreach ($r in $rows) {
        $row = $r.Row
        $guid = $r.Guid

        $mdName = "$guid.md"
        $xlfName = "$guid.42cee4cdfc434bab753d64562e2f8705da84442b.$Lang.xlf"

        $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d575ad2a16b5c33db1b7539e4b96a2a654c4ad9e/e2e/$mdName"
        $handbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$HandoffCommit/ol-handback/OpenLocalizationTestOrg/oltest.$Lang/ci/hb/$xlfName"

        # Status -> handed back, in sync with en-US
        $ws.Range("C$row").Value = $StatusText

        # Latest Target File (new column F): same source .md, linked back
        # to the source repo like column A / B.
        $fCell = $ws.Range("F$row")
        $fCell.Value = $mdName
        $ws.Hyperlinks.Add($fCell, $mdUrl, "", "", $mdName)
        $fCell.Style = "HyperLink"

        # Latest Handback File (new column G): the handed-back xlf.
        $gCell = $ws.Range("G$row")
        $gCell.Value = $xlfName
        $ws.Hyperlinks.Add($gCell, $handbackUrl, "", "", $xlfName)
        $gCell.Style = "HyperLink"

        # Latest Handback DateTime (column H): was the zero-date
        # placeholder, now the real handback timestamp.
        $ws.Range("H$row").Value = $HandbackDateTime
    }
}

Update-LangSheet $wb "zh-cn" "zh-cn" "2016-03-21 00:52:47" "821ddd3a7df75ea26a1890bd27fc56443411af55" $statusText
Update-LangSheet $wb "de-de" "de-de" "2016-03-21 00:52:53" "dd036ff46da9d4372695b6e9c0ea92502a7338f9" $statusText
